$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# New weekly data rows appended to the Main sheet (rows 127-136).
# Each hashtable maps column letter -> value; missing columns stay blank.
$rows = @(
    @{ rowNum = 127; A = 887; B = 311;  C = 5;     D = 17514; E = 17839; F = 16576;            H = 17125; I = 15863;            K = 20399; L = 22399; M = 3475;  N = 4;  O = 1;  P = -2; Q = -1;  R = 4;  S = 0;  T = 0 },
    @{ rowNum = 128; A = 888; B = 2111; C = 32;    D = 18981; E = 17290; F = 17270; G = 16559; H = 16605; I = 18569; J = 15594; K = 13028;            M = 16227; N = 2;  O = 5;  P = -2; Q = 1;   R = -1; S = 2;  T = -1 },
    @{ rowNum = 129; A = 889; B = 2111; C = 2111;  D = 16866; E = 16556; F = 16785; G = 15488; H = 15752; I = 16891; J = 16232; K = 15165;            M = 15646; N = 0;  O = -2; P = 1;  Q = -4;  R = 0;  S = 2;  T = 3 },
    @{ rowNum = 130; A = 890; B = 2111; C = 32;    D = 18528; E = 17161; F = 17101; G = 16939; H = 16499; I = 19444; J = 15136; K = 21108; L = 12868; M = 16422; N = 1;  O = -2; P = -1; Q = -8;  R = 11; S = 3;  T = -4 },
    @{ rowNum = 131; A = 891; B = 2111; C = 221;   D = 16497; E = 16104; F = 16120; G = 15538; H = 15305;            J = 18709; K = 15525; L = 15676;            N = 0;  O = 4;  P = 1;  Q = 2;   R = 1;  S = 1;  T = 2 },
    @{ rowNum = 132; A = 892; B = 2111; C = 5;     D = 16456; E = 15517; F = 17018; G = 16104; H = 15123; I = 20815; J = 10000; K = 16530; L = 16903; M = 10891; N = -8; O = 6;  P = 2;  Q = -10; R = 6;  S = 1;  T = 3 },
    @{ rowNum = 133; A = 893; B = 2111; C = 5;     D = 16989; E = 16662; F = 16769; G = 16400; H = 16400; I = 17989; J = 19879; K = 15044; L = 16507; M = 13873; N = -5; O = 17; P = 1;  Q = -3;  R = 4;  S = -7; T = -7 },
    @{ rowNum = 134; A = 894; B = 5;    C = 2111;  D = 23959; E = 14999; F = 22695; G = 16750; H = 23331; I = 22263;            K = 21617; L = 23149;            N = 3;  O = 7;  P = 0;  Q = -1;  R = 0;  S = -4; T = 4 },
    @{ rowNum = 135; A = 895; B = 5;    C = 221;   D = 23836; E = 22595; F = 16454; G = 22994;            I = 21505; J = 23400; K = 19999; L = 24783; M = 22248; N = -3; O = 3;  P = 2;  Q = -5;  R = -4; S = -3; T = 8 },
    @{ rowNum = 136; A = 896; B = 41;   C = 11111; D = 22471; E = 23520;            G = 16147; H = 19882; I = 21183; J = 20960; K = 19452; L = 20005; M = 19999; N = -4; O = 1;  P = 0;  Q = -3;  R = -2; S = -4; T = 10 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rows) {
    foreach ($col in $cols) {
        if ($row.ContainsKey($col)) {
            $ws.Range("$col$($row.rowNum)").Value = $row[$col]
        }
    }
}

# Sheet view updates: frozen pane top-left cell and the active selection.
$ws.Application.ActiveWindow.ScrollRow = 118
$ws.Range("P136").Select()
